# Update latest output (run 140)
# Applies new optimisation-result values to the "Schedule" and "Detailed" sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet "Schedule": recalculated cost / unit-cost figures ----
$schedule = $wb.Worksheets.Item("Schedule")

$schedule.Range("E2").Value = 506.74358475
$schedule.Range("F2").Value = 33.51478735119048

$schedule.Range("E3").Value = -239.3418885
$schedule.Range("F3").Value = -7.914744990079364

$schedule.Range("E4").Value = 497.6241465
$schedule.Range("F4").Value = 32.91164990079366

# ---- Sheet "Detailed": refreshed price series + historical/forecast split ----
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B9").Value = 65.64212999999999
$detailed.Range("B10").Value = 68.87112999999999

$detailed.Range("C11").Value = "historical"

$detailed.Range("B12").Value = 96.67227
$detailed.Range("C12").Value = "historical"

$detailed.Range("C13").Value = "historical"

$detailed.Range("B15").Value = 78
$detailed.Range("B16").Value = 57.03041
$detailed.Range("B17").Value = 12.39286
$detailed.Range("B18").Value = 0

$detailed.Range("B20").Value = -6.73725
$detailed.Range("B21").Value = -7.34167
$detailed.Range("B22").Value = -8.422269999999999
$detailed.Range("B23").Value = -10.2958
$detailed.Range("B24").Value = -15.15889
$detailed.Range("B25").Value = -15.60123
$detailed.Range("B26").Value = -15.56494
$detailed.Range("B27").Value = -21.65844
$detailed.Range("B28").Value = -21.24892
$detailed.Range("B29").Value = -22.48391
$detailed.Range("B30").Value = -25.26417
$detailed.Range("B31").Value = -23.8909
$detailed.Range("B32").Value = -23.5
$detailed.Range("B33").Value = -22.56734
$detailed.Range("B34").Value = -6.75345

$detailed.Range("B36").Value = 36.06
$detailed.Range("B37").Value = 47.6595
$detailed.Range("B38").Value = 47.92523
$detailed.Range("B39").Value = 57.31
$detailed.Range("B40").Value = 71.21621
$detailed.Range("B41").Value = 73.19
$detailed.Range("B42").Value = 77.94

$detailed.Range("B44").Value = 57.31
$detailed.Range("B45").Value = 59.86017

$detailed.Range("B47").Value = 63.18337
